# OrdLine.xlsx correction: two OrdNo typos fixed in column A.
#   Row 4  : "O1233131" -> "O1231231"  (matches the O1231231 order already on row 3)
#   Row 23 : "O333222"  -> "O3331222"  (matches the O3331222 order already on row 22)
#   Row 24 : "O333222"  -> "O3331222"  (same order as row 23, second product line)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "O1231231"
$ws.Range("A23").Value = "O3331222"
$ws.Range("A24").Value = "O3331222"

# Let Excel re-measure the (word-wrapped) row heights for the rows whose
# text just changed length, matching what a live edit in Excel would do.
$ws.Rows("23:24").AutoFit()

# Leave the selection where the author left it before saving.
$ws.Range("R9").Select()
